# Rename the worksheet from "Output results" to "output results".
# This also updates the _FilterDatabase defined name reference automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "output results"

# Remove the wrap-text alignment from F4:AK4 so these cells share the
# same style as the equivalent cells in rows 3 and 5.
$ws.Range("F4:AK4").WrapText = $false
